$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E4").Value = "值我"
$ws.Range("E5").Value = "聖祖仁旬大慶節欽奉"
$ws.Range("E2").Value = "敕大乾海含弘至德四位上"
